$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '69.969.13'
$ws.Range("E2").Value = '  -0.58%  '
$ws.Range("D3").Value = '3.535.05'
$ws.Range("E3").Value = '  -0.97%  '
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '604.20'
$ws.Range("E5").Value = '  +2.38%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '184.57'
$ws.Range("E6").Value = '  -1.34%  '
$ws.Range("D7").Value = '3.528.87'
$ws.Range("E7").Value = '  -0.80%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.614'
$ws.Range("E8").Value = '  -1.18%  '
$ws.Range("E9").Value = '  -0.10%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.212'
$ws.Range("E10").Value = '  +6.15%  '
$ws.Range("E11").Value = '  -1.24%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '53.49'
$ws.Range("E12").Value = '  -2.30%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000307'
$ws.Range("E13").Value = '  -0.35%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '9.42'
$ws.Range("E14").Value = '  -1.80%  '
$ws.Range("D15").Value = '4.110.06'
$ws.Range("E15").Value = '  -0.64%  '
$ws.Range("D16").Value = '69.991.51'
$ws.Range("E16").Value = '  -0.53%  '
$ws.Range("D17").Value = '3.568.88'
$ws.Range("E17").Value = '  +0.25%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '18.86'
$ws.Range("E18").Value = '  -3.17%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.59'
$ws.Range("E19").Value = '  +0.87%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '576.78'
$ws.Range("E20").Value = '  +4.82%  '
$ws.Range("E21").Value = '  -0.03%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.986'
$ws.Range("E22").Value = '  -3.46%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '17.29'
$ws.Range("E23").Value = '  -4.01%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '4.67'
$ws.Range("E24").Value = '  -0.04%  '
$ws.Range("E25").Value = '  -2.13%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '94.15'
$ws.Range("E26").Value = '  -2.07%  '
$ws.Range("E27").Value = '  -2.31%  '
$ws.Range("E28").Value = '  -4.81%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.29'
$ws.Range("E29").Value = '  +1.26%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '31.98'
$ws.Range("E30").Value = '  -0.82%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.99'
$ws.Range("E31").Value = '  -4.93%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '12.17'
$ws.Range("E32").Value = '  -3.32%  '
$ws.Range("E33").Value = '  -1.34%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '63.28'
$ws.Range("E34").Value = '  -2.96%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.64'
$ws.Range("E35").Value = '  +17.69%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '3.27'
$ws.Range("E36").Value = '  +1.03%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '530.31'
$ws.Range("E37").Value = '  -3.92%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.401'
$ws.Range("E38").Value = '  -4.27%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.00'
$ws.Range("E39").Value = '  +0.04%  '
$ws.Range("E40").Value = '  -3.99%  '
$ws.Range("D41").Value = '0.0₃0777'
$ws.Range("E41").Value = '  +0.99%  '
$ws.Range("D42").Value = '3.524.71'
$ws.Range("E42").Value = '  +4.41%  '
$ws.Range("E43").Value = '  +3.13%  '
$ws.Range("E44").Value = '  +0.11%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0454'
$ws.Range("E45").Value = '  +1.09%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.43'
$ws.Range("E46").Value = '  -4.35%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.89'
$ws.Range("E47").Value = '  -3.16%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.139'
$ws.Range("E48").Value = '  +2.06%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '9.13'
$ws.Range("E49").Value = '  -0.88%  '
$ws.Range("E50").Value = '  +0.24%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.41'
$ws.Range("E51").Value = '  -5.49%  '
